$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Bump every student's "B" column number (rows 2-25) up by 7.
for ($r = 2; $r -le 25; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value2 = $cell.Value2 + 7
}

# Scroll the frozen pane back to the top (topLeftCell A17 -> A2) and
# move the selection to E29, matching the saved view state.
$ws.Range("A2").Select()
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("E29").Select()
